# ADD results from server
# Updates the cached result values on each year sheet ("2025".."2050")
# with the refreshed numbers produced by the server run.

$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{
        "B2" = 973.9537847600009
        "E2" = 28982.37596598056
        "I2" = 16175.28135478
        "L2" = 48524.529503538
        "M2" = 10590.587968015
        "N2" = 7155.07579047334
        "O2" = 6980.325566461758
    }
    "2030" = @{
        "B2" = 5712.560177842886
        "E2" = 56106.05588781912
        "I2" = 44217.8984721661
        "L2" = 66966.57749858923
        "M2" = 21984.28023276101
        "N2" = 10592.74688452318
        "O2" = 12062.16651258332
    }
    "2035" = @{
        "A2" = 2861.961401238371
        "B2" = 8026.889663087295
        "E2" = 67297.73995507321
        "I2" = 59256.42575923612
        "L2" = 66966.57749858923
        "M2" = 25464.6214365565
        "N2" = 15131.83576261564
        "O2" = 14760.26691511388
    }
    "2040" = @{
        "A2" = 2861.961401238371
        "B2" = 8026.889663087295
        "E2" = 67297.73995507321
        "I2" = 59256.42575923612
        "L2" = 66966.57749858923
        "M2" = 25464.6214365565
        "N2" = 15236.24445882674
        "O2" = 14760.26691511388
    }
    "2045" = @{
        "A2" = 6302.873118834019
        "B2" = 8026.889663087295
        "E2" = 67297.73995507321
        "I2" = 59256.42575923612
        "L2" = 66966.57749858923
        "M2" = 25464.6214365565
        "N2" = 15771.75564340065
        "O2" = 17095.73290146254
    }
    "2050" = @{
        "A2" = 6302.873118834019
        "B2" = 8026.889663087295
        "E2" = 67297.73995507321
        "I2" = 59256.42575923612
        "L2" = 66966.57749858923
        "M2" = 25464.6214365565
        "N2" = 15771.75564340065
        "O2" = 17095.73290146254
    }
}

foreach ($sheetName in $updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellValues = $updates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range([string]$cellRef).Value = $cellValues[$cellRef]
    }
}
